# Update 16S tree with Sanger contig sequence
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: "Sanger sequence" -> "Sanger_ID"
$ws.Range("H1").Value = "Sanger_ID"

# Record the Sanger sequencing IDs (1-5) for the five successfully
# extracted/classified isolates, and the classification result for C05.
$ws.Range("H14").Value = 1
$ws.Range("H14").Style = "Good"

$ws.Range("H18").Value = 2
$ws.Range("H18").Style = "Good"
$ws.Range("I18").Value = "Contig_2: Ramlibacter"
$ws.Range("I18").Style = "Good"

$ws.Range("H23").Value = 3
$ws.Range("H23").Style = "Good"

$ws.Range("H24").Value = 4
$ws.Range("H24").Style = "Good"

$ws.Range("H25").Value = 5
$ws.Range("H25").Style = "Good"

# Widen column I to fit the new classification text
$ws.Columns.Item(9).ColumnWidth = 20.5703125

# Move / extend the active selection to the newly filled-in row
$excel.Goto($ws.Range("C18:H18"))
